$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

# Append two new rows of "people" data (new singleton-style entries),
# following the existing idPers / nomPers pattern already in the sheet.
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "@mickaelbaron"
$ws.Cells.Item(23, 3).NumberFormat = "General"

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "mikiu"
$ws.Cells.Item(24, 3).NumberFormat = "General"
